$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells that need to be marked absent ("A") that were previously blank.
$cells = @(
    "V7",
    "W10",
    "W11",
    "W13",
    "V14",
    "W15",
    "V17",
    "W17",
    "W19",
    "W25",
    "V26",
    "W28",
    "W29",
    "V32",
    "W33",
    "W34",
    "W36",
    "V37",
    "W38",
    "W39",
    "W41",
    "W42",
    "V51",
    "W54",
    "W60",
    "W64",
    "W71",
    "V78"
)

foreach ($cell in $cells) {
    $ws.Range($cell).Value = "A"
    $ws.Range($cell).Orientation = 0
}
